# Auto-update draw results: append the 2025-10-22 Pick 3 result as a new
# row right after the last existing row of data.
#
# The sheet stores every column as literal text (dates, phase codes and the
# result string all look numeric/date-like), so each value is written with a
# leading apostrophe to force Excel to keep it as text instead of silently
# re-interpreting it as a date serial or a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 36

$ws.Range("A$newRow").Value = "'2025-10-22"
$ws.Range("B$newRow").Value = "Pick 3"
$ws.Range("C$newRow").Value = "'251022"
$ws.Range("D$newRow").Value = "'8-5-2"
$ws.Range("E$newRow").Value = "'2025-10-22T21:37:35.596+04:00"
